$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "25.983.40"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "1.641.07"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.58%  "
$ws.Range("D5").Value = "215.97"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Value = "1.01"
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").Value = "0.0637"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "19.61"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "1.868.24"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").Value = "1.665.79"
$ws.Range("E14").Value = "  +2.17%  "
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "25.984.95"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").Value = "1.01"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").Value = "192.98"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").Value = "9.92"
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D23").Value = "6.28"
$ws.Range("E23").Value = "  -0.72%  "
$ws.Range("E24").Value = "  +4.33%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("D27").Value = "143.61"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("D28").Value = "6.88"
$ws.Range("E28").Value = "  -1.54%  "
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  -0.83%  "
$ws.Range("D31").Value = "0.0502"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("E32").Value = "  -1.73%  "
$ws.Range("D33").Value = "3.24"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("E34").Value = "  -5.12%  "
$ws.Range("E35").Value = "  +1.95%  "
$ws.Range("E36").Value = "  -1.50%  "
$ws.Range("D37").Value = "1.135.13"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "0.543"
$ws.Range("E38").Value = "  -2.44%  "
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").Value = "1.01"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("E42").Value = "  -1.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.30"
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").Value = "1.777.54"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("E46").Value = "  +1.92%  "
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0530"
$ws.Range("E48").Value = "  +2.36%  "
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").Value = "0.416"
$ws.Range("E51").Value = "  -0.43%  "
